# Auto-generated Excel COM-interop script
# Applies numeric corrections to the H:N profit-calculation columns
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# per the scheduled-runner profit refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 399.75
$ws.Range("I19").Value = 349.5
$ws.Range("K19").Value = 349.5
$ws.Range("M19").Value = -174.5
$ws.Range("H33").Value = 135.83333
$ws.Range("I33").Value = 74.666664
$ws.Range("J33").Value = 197
$ws.Range("K33").Value = 74.666664
$ws.Range("L33").Value = 197
$ws.Range("M33").Value = 154.333336
$ws.Range("N33").Value = -655
$ws.Range("H53").Value = 371.81818
$ws.Range("J53").Value = 132.5
$ws.Range("L53").Value = 132.5
$ws.Range("N53").Value = -1406.5
$ws.Range("H70").Value = 5999.625
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 6999.5
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 20998.5
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -21538.5
$ws.Range("H73").Value = 5999.625
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 6999.5
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 20998.5
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -22870.5
$ws.Range("H96").Value = 214.42857
$ws.Range("I96").Value = 260.4
$ws.Range("K96").Value = 781.1999999999999
$ws.Range("M96").Value = 591.8000000000001
$ws.Range("H116").Value = 4766
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 8584.173000000001
$ws.Range("I132").Value = 8553.370000000001
$ws.Range("K132").Value = 25660.11
$ws.Range("M132").Value = -23130.11
$ws.Range("H135").Value = 643.1429000000001
$ws.Range("J135").Value = 320
$ws.Range("L135").Value = 2880
$ws.Range("N135").Value = -7950
$ws.Range("H137").Value = 1921.8235
$ws.Range("I137").Value = 891.8
$ws.Range("K137").Value = 2675.4
$ws.Range("M137").Value = -125.3999999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1767.375
$ws.Range("I61").Value = 1328.7693
$ws.Range("K61").Value = 1328.7693
$ws.Range("M61").Value = -1116.7693
$ws.Range("H74").Value = 2526
$ws.Range("I74").Value = 1748.2142
$ws.Range("K74").Value = 1748.2142
$ws.Range("M74").Value = -874.2141999999999
$ws.Range("H77").Value = 2526
$ws.Range("I77").Value = 1748.2142
$ws.Range("K77").Value = 8741.071
$ws.Range("M77").Value = -4373.071
$ws.Range("H132").Value = 1780.7894
$ws.Range("I132").Value = 1780.7894
$ws.Range("K132").Value = 5342.3682
$ws.Range("M132").Value = -2812.3682
$ws.Range("H136").Value = 1767.375
$ws.Range("I136").Value = 1328.7693
$ws.Range("K136").Value = 3986.3079
$ws.Range("M136").Value = -1436.3079
$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 906.25
$ws.Range("I94").Value = 906.25
$ws.Range("K94").Value = 906.25
$ws.Range("M94").Value = -455.25
$ws.Range("H99").Value = 4006.6667
$ws.Range("I99").Value = 4006.6667
$ws.Range("K99").Value = 4006.6667
$ws.Range("M99").Value = -2508.6667
$ws.Range("H134").Value = 2195.875
$ws.Range("I134").Value = 942.8333
$ws.Range("J134").Value = 5955
$ws.Range("K134").Value = 2828.4999
$ws.Range("L134").Value = 17865
$ws.Range("M134").Value = -293.4998999999998
$ws.Range("N134").Value = -22935

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7636.2
$ws.Range("I31").Value = 2532.25
$ws.Range("J31").Value = 9492.182000000001
$ws.Range("K31").Value = 2532.25
$ws.Range("L31").Value = 9492.182000000001
$ws.Range("M31").Value = -2237.25
$ws.Range("N31").Value = -10082.182
$ws.Range("H34").Value = 7636.2
$ws.Range("I34").Value = 2532.25
$ws.Range("J34").Value = 9492.182000000001
$ws.Range("K34").Value = 2532.25
$ws.Range("L34").Value = 9492.182000000001
$ws.Range("M34").Value = -2330.25
$ws.Range("N34").Value = -9896.182000000001
$ws.Range("H58").Value = 2593.7058
$ws.Range("I58").Value = 1309.5834
$ws.Range("K58").Value = 1309.5834
$ws.Range("M58").Value = -1106.5834
$ws.Range("H105").Value = 3005.5557
$ws.Range("I105").Value = 3005.5557
$ws.Range("K105").Value = 3005.5557
$ws.Range("M105").Value = -1258.5557
$ws.Range("H132").Value = 2040.2
$ws.Range("I132").Value = 2181.5
$ws.Range("J132").Value = 1475
$ws.Range("K132").Value = 6544.5
$ws.Range("L132").Value = 4425
$ws.Range("M132").Value = -4014.5
$ws.Range("N132").Value = -9485
$ws.Range("H134").Value = 4740.2856
$ws.Range("I134").Value = 2833.25
$ws.Range("K134").Value = 8499.75
$ws.Range("M134").Value = -5964.75
$ws.Range("H136").Value = 2593.7058
$ws.Range("I136").Value = 1309.5834
$ws.Range("K136").Value = 3928.7502
$ws.Range("M136").Value = -1378.7502

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4999.8
$ws.Range("J19").Value = 4999.8
$ws.Range("L19").Value = 14999.4
$ws.Range("N19").Value = -15347.4
$ws.Range("H38").Value = 459.9524
$ws.Range("I38").Value = 436.55554
$ws.Range("J38").Value = 600.3333
$ws.Range("K38").Value = 1309.66662
$ws.Range("L38").Value = 1800.9999
$ws.Range("M38").Value = -962.66662
$ws.Range("N38").Value = -2494.9999
$ws.Range("H55").Value = 7053.75
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 7053.75
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21161.25
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -21515.25
$ws.Range("H113").Value = 1349.0714
$ws.Range("I113").Value = 929.625
$ws.Range("K113").Value = 2788.875
$ws.Range("M113").Value = -618.875
$ws.Range("H122").Value = 749.6667
$ws.Range("I122").Value = 749
$ws.Range("K122").Value = 6741
$ws.Range("M122").Value = -4291
$ws.Range("H131").Value = 982.3333
$ws.Range("I131").Value = 982.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2946.9999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2093.0001
$ws.Range("N131").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820
$ws.Range("H126").Value = 3303.5557
$ws.Range("I126").Value = 2480.4
$ws.Range("J126").Value = 4332.5
$ws.Range("K126").Value = 7441.200000000001
$ws.Range("L126").Value = 12997.5
$ws.Range("M126").Value = -4971.200000000001
$ws.Range("N126").Value = -17937.5
$ws.Range("H132").Value = 2456.375
$ws.Range("I132").Value = 2164.5715
$ws.Range("K132").Value = 6493.7145
$ws.Range("M132").Value = -3963.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 515000
$ws.Range("I69").Value = 30000
$ws.Range("J69").Value = 1000000
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 1000000
$ws.Range("M69").Value = -29189
$ws.Range("N69").Value = -1001622
$ws.Range("H72").Value = 515000
$ws.Range("I72").Value = 30000
$ws.Range("J72").Value = 1000000
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 3000000
$ws.Range("M72").Value = -85944
$ws.Range("N72").Value = -3008112
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H93").Value = 1340.7273
$ws.Range("I93").Value = 1507.1428
$ws.Range("J93").Value = 1049.5
$ws.Range("K93").Value = 1507.1428
$ws.Range("L93").Value = 1049.5
$ws.Range("M93").Value = -259.1428000000001
$ws.Range("N93").Value = -3545.5
$ws.Range("H122").Value = 2784.1538
$ws.Range("I122").Value = 2724.5
$ws.Range("K122").Value = 8173.5
$ws.Range("M122").Value = -5723.5
$ws.Range("H132").Value = 4395.25
$ws.Range("I132").Value = 4326.7144
$ws.Range("K132").Value = 12980.1432
$ws.Range("M132").Value = -10450.1432
$ws.Range("H136").Value = 1623.3334
$ws.Range("I136").Value = 1648
$ws.Range("K136").Value = 4944
$ws.Range("M136").Value = -2394

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8181.8184
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 8181.8184
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H122").Value = 1254.7368
$ws.Range("I122").Value = 1121.3334
$ws.Range("K122").Value = 3364.0002
$ws.Range("M122").Value = -914.0001999999999
$ws.Range("H132").Value = 1826.72
$ws.Range("I132").Value = 1736.1666
$ws.Range("K132").Value = 5208.4998
$ws.Range("M132").Value = -2678.4998
$ws.Range("H136").Value = 3495.6667
$ws.Range("I136").Value = 2802.739
$ws.Range("K136").Value = 8408.217000000001
$ws.Range("M136").Value = -5858.217000000001
$ws.Range("H138").Value = 101750
$ws.Range("J138").Value = 101750
$ws.Range("L138").Value = 101750
$ws.Range("N138").Value = -112030
$ws.Range("H141").Value = 171176.62
$ws.Range("J141").Value = 124201.86
$ws.Range("L141").Value = 124201.86
$ws.Range("N141").Value = -134561.86

